# Worker List.xlsx update
# - Fix scanner ID typo for Justin Blake (row 19): 63F447B -> 637F447B
# - Fill in missing scanner assignments for Korey Sugar (row 21),
#   Mark Zhukov (row 23) and Saad Shamsaldeen (row 30)
# - Update the active cell selection left on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Korey Sugar was missing his S/N and ID values
$ws.Range("C21").Value = "1AF41C46"
$ws.Range("D21").Value = "N521D5060014"

# Mark Zhukov was missing his S/N value (placeholder cell was centered);
# match the left-aligned formatting used by the other ID cells in column C
$ws.Range("C23").HorizontalAlignment = -4131
$ws.Range("C23").Value = "FB2B1B1E"
$ws.Range("D23").Value = "N521D5060017"

# Saad Shamsaldeen was missing his S/N and ID values
$ws.Range("C30").Value = "D1CECBAB"
$ws.Range("D30").Value = "N521D5060029"

# Correct the typo'd serial number for Justin Blake
$ws.Range("C19").Value = "637F447B"

# Leave the cursor positioned where the author left off editing
$ws.Range("H15").Select() | Out-Null
